$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the missing "SCALE D" readings (column F) for rows 14-19 ---
$ws.Range("F14").Value = 10498
$ws.Range("F15").Value = 10495
$ws.Range("F16").Value = 10484
$ws.Range("F17").Value = 10474
$ws.Range("F18").Value = 10507
$ws.Range("F19").Value = 10428

# --- Correct the "Gnome Weight" values (column G) for rows 14-19 ---
$ws.Range("G14").Value = 22.7615
$ws.Range("G15").Value = 22.7615
$ws.Range("G16").Value = 22.7615
$ws.Range("G17").Value = 22.7615
$ws.Range("G18").Value = 22.7615
$ws.Range("G19").Value = 22.7615

# --- Remove the stray manual-offset value that was sitting in H21 ---
$ws.Range("H21").ClearContents()

# --- Update the window/view position so the sheet scrolls down to show
#     the newly completed SCALE D table, with H20 as the active cell ---
$ws.Range("H20").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
